# Case_1_160 res_line/pl_mw.xlsx update: "case with 380 kV done"
# Updates line-power-flow results (columns B,D,E,F,G,H,I,J,L) for rows 2-25
# on Sheet1 to the newly recomputed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.76820097031981
$ws.Range("D2").Value = 0.145921414004988
$ws.Range("E2").Value = 0.1814845868011474
$ws.Range("F2").Value = 1.742266804659053
$ws.Range("G2").Value = 1.274910148976801
$ws.Range("H2").Value = 1.165069515157768
$ws.Range("I2").Value = 0.9147059153448183
$ws.Range("J2").Value = 0.256761826735584
$ws.Range("L2").Value = 0.9989631696511765
# Row 3
$ws.Range("B3").Value = 1.635495077113944
$ws.Range("D3").Value = 0.1396554260951746
$ws.Range("E3").Value = 0.1716639927178107
$ws.Range("F3").Value = 1.747959765855725
$ws.Range("G3").Value = 1.260638274770471
$ws.Range("H3").Value = 1.166689257361639
$ws.Range("I3").Value = 0.9369197368218689
$ws.Range("J3").Value = 0.2411312169825663
$ws.Range("L3").Value = 0.911455373729666
# Row 4
$ws.Range("B4").Value = 1.554045017029011
$ws.Range("D4").Value = 0.1357749990589383
$ws.Range("E4").Value = 0.1656119241487843
$ws.Range("F4").Value = 1.753005444373173
$ws.Range("G4").Value = 1.253198973076877
$ws.Range("H4").Value = 1.168604239977711
$ws.Range("I4").Value = 0.9513939186398197
$ws.Range("J4").Value = 0.2315154996156537
$ws.Range("L4").Value = 0.8577648477051412
# Row 5
$ws.Range("B5").Value = 1.520863088655176
$ws.Range("D5").Value = 0.1341853262824344
$ws.Range("E5").Value = 0.163140192663576
$ws.Range("F5").Value = 1.755449729365822
$ws.Range("G5").Value = 1.250498024972856
$ws.Range("H5").Value = 1.16961519780385
$ws.Range("I5").Value = 0.957501761153889
$ws.Range("J5").Value = 0.2275926394208199
$ws.Range("L5").Value = 0.8358961780160143
# Row 6
$ws.Range("B6").Value = 1.515353884890089
$ws.Range("D6").Value = 0.1339208556097162
$ws.Range("E6").Value = 0.1627294369136081
$ws.Range("F6").Value = 1.75587899424724
$ws.Range("G6").Value = 1.250069438797823
$ws.Range("H6").Value = 1.169796968137206
$ws.Range("I6").Value = 0.9585286045259718
$ws.Range("J6").Value = 0.2269409949871886
$ws.Range("L6").Value = 0.8322655666956393
# Row 7
$ws.Range("B7").Value = 1.553597472777881
$ws.Range("D7").Value = 0.1357535941135382
$ws.Range("E7").Value = 0.1655786114832907
$ws.Range("F7").Value = 1.7530368395647
$ws.Range("G7").Value = 1.25316121139619
$ws.Range("H7").Value = 1.16861694163525
$ws.Range("I7").Value = 0.9514754434649255
$ws.Range("J7").Value = 0.2314626119178484
$ws.Range("L7").Value = 0.8574698750901177
# Row 8
$ws.Range("B8").Value = 1.722438270652901
$ws.Range("D8").Value = 0.143767737332098
$ws.Range("E8").Value = 0.1781031148335188
$ws.Range("F8").Value = 1.743906984385575
$ws.Range("G8").Value = 1.269713126299536
$ws.Range("H8").Value = 1.165436434291394
$ws.Range("I8").Value = 0.9221917234627934
$ws.Range("J8").Value = 0.2513763761489543
$ws.Range("L8").Value = 0.9687826405036049
# Row 9
$ws.Range("B9").Value = 2.053735487917493
$ws.Range("D9").Value = 0.1592238383876463
$ws.Range("E9").Value = 0.2024838999032426
$ws.Range("F9").Value = 1.738381172878903
$ws.Range("G9").Value = 1.312778510812564
$ws.Range("H9").Value = 1.166543806661679
$ws.Range("I9").Value = 0.8714086726705617
$ws.Range("J9").Value = 0.2902712355528507
$ws.Range("L9").Value = 1.187362822326122
# Row 10
$ws.Range("B10").Value = 2.297219529737333
$ws.Range("D10").Value = 0.170426282177047
$ws.Range("E10").Value = 0.2202838334427852
$ws.Range("F10").Value = 1.741978332391668
$ws.Range("G10").Value = 1.351035113415662
$ws.Range("H10").Value = 1.171893935994746
$ws.Range("I10").Value = 0.8381722413346662
$ws.Range("J10").Value = 0.3187425609097119
$ws.Range("L10").Value = 1.348128157536109
# Row 11
$ws.Range("B11").Value = 2.407997545043941
$ws.Range("D11").Value = 0.1754904994927102
$ws.Range("E11").Value = 0.2283565462464878
$ws.Range("F11").Value = 1.745301736484251
$ws.Range("G11").Value = 1.369908606226062
$ws.Range("H11").Value = 1.175325962965275
$ws.Range("I11").Value = 0.8239421289358848
$ws.Range("J11").Value = 0.3316703985678515
$ws.Range("L11").Value = 1.421302686146475
# Row 12
$ws.Range("B12").Value = 2.449947540577227
$ws.Range("D12").Value = 0.1774036858385131
$ws.Range("E12").Value = 0.231409864107178
$ws.Range("F12").Value = 1.746804761538726
$ws.Range("G12").Value = 1.377269528107661
$ws.Range("H12").Value = 1.176770159585089
$ws.Range("I12").Value = 0.8186819329785475
$ws.Range("J12").Value = 0.3365621927207627
$ws.Range("L12").Value = 1.449017645412823
# Row 13
$ws.Range("B13").Value = 2.440912848121911
$ws.Range("D13").Value = 0.1769918471961205
$ws.Range("E13").Value = 0.2307524412388062
$ws.Range("F13").Value = 1.746470152036096
$ws.Range("G13").Value = 1.375674666810397
$ws.Range("H13").Value = 1.176452680215903
$ws.Range("I13").Value = 0.8198090889725922
$ws.Range("J13").Value = 0.3355088258299332
$ws.Range("L13").Value = 1.443048504574335
# Row 14
$ws.Range("B14").Value = 2.411448789627912
$ws.Range("D14").Value = 0.1756479888100415
$ws.Range("E14").Value = 0.228607818183292
$ws.Range("F14").Value = 1.745420479528391
$ws.Range("G14").Value = 1.370509891631883
$ws.Range("H14").Value = 1.175441874773355
$ws.Range("I14").Value = 0.8235067921203232
$ws.Range("J14").Value = 0.3320729249808778
$ws.Range("L14").Value = 1.423582708702952
# Row 15
$ws.Range("B15").Value = 2.393401258387428
$ws.Range("D15").Value = 0.1748242493856935
$ws.Range("E15").Value = 0.2272936969464396
$ws.Range("F15").Value = 1.74480942658208
$ws.Range("G15").Value = 1.367374253544369
$ws.Range("H15").Value = 1.174841583956464
$ws.Range("I15").Value = 0.8257884852814463
$ws.Range("J15").Value = 0.3299678474042196
$ws.Range("L15").Value = 1.411660025828439
# Row 16
$ws.Range("B16").Value = 2.289980114743742
$ws.Range("D16").Value = 0.1700946877546698
$ws.Range("E16").Value = 0.2197557588039416
$ws.Range("F16").Value = 1.741795252296598
$ws.Range("G16").Value = 1.349831487392265
$ws.Range("H16").Value = 1.171689818691277
$ws.Range("I16").Value = 0.8391201552337257
$ws.Range("J16").Value = 0.3178971928888075
$ws.Range("L16").Value = 1.343346812432515
# Row 17
$ws.Range("B17").Value = 2.226537568037884
$ws.Range("D17").Value = 0.1671851403671951
$ws.Range("E17").Value = 0.2151251061141792
$ws.Range("F17").Value = 1.740379613676225
$ws.Range("G17").Value = 1.339447940739717
$ws.Range("H17").Value = 1.170012704287956
$ws.Range("I17").Value = 0.8475268696763294
$ws.Range("J17").Value = 0.3104859269655122
$ws.Range("L17").Value = 1.301449061973074
# Row 18
$ws.Range("B18").Value = 2.19004867791125
$ws.Range("D18").Value = 0.1655086462464794
$ws.Range("E18").Value = 0.2124593712791523
$ws.Range("F18").Value = 1.739724051607624
$ws.Range("G18").Value = 1.333613847119892
$ws.Range("H18").Value = 1.169141968761323
$ws.Range("I18").Value = 0.8524458614018364
$ws.Range("J18").Value = 0.3062209273147971
$ws.Range("L18").Value = 1.277354532327649
# Row 19
$ws.Range("B19").Value = 2.177694486028031
$ws.Range("D19").Value = 0.1649404973691873
$ws.Range("E19").Value = 0.2115564066575999
$ws.Range("F19").Value = 1.73952928368486
$ws.Range("G19").Value = 1.331662191497003
$ws.Range("H19").Value = 1.168863248558694
$ws.Range("I19").Value = 0.8541257013572405
$ws.Range("J19").Value = 0.3047764957548935
$ws.Range("L19").Value = 1.269197236247976
# Row 20
$ws.Range("B20").Value = 2.23329098259876
$ws.Range("D20").Value = 0.1674951766935777
$ws.Range("E20").Value = 0.2156182864602485
$ws.Range("F20").Value = 1.740513874548768
$ws.Range("G20").Value = 1.340538960824233
$ws.Range("H20").Value = 1.170181510720596
$ws.Range("I20").Value = 0.8466232968379348
$ws.Range("J20").Value = 0.3112751021366336
$ws.Range("L20").Value = 1.305908745023999
# Row 21
$ws.Range("B21").Value = 2.420103088938333
$ws.Range("D21").Value = 0.1760428348992207
$ws.Range("E21").Value = 0.2292378457045459
$ws.Range("F21").Value = 1.745722142021734
$ws.Range("G21").Value = 1.372021085002871
$ws.Range("H21").Value = 1.175734841280985
$ws.Range("I21").Value = 0.8224171967935874
$ws.Range("J21").Value = 0.3330822347293747
$ws.Range("L21").Value = 1.429300144343358
# Row 22
$ws.Range("B22").Value = 2.542199389127461
$ws.Range("D22").Value = 0.1816028877211977
$ws.Range("E22").Value = 0.2381177301521973
$ws.Range("F22").Value = 1.750552225270027
$ws.Range("G22").Value = 1.393844652024882
$ws.Range("H22").Value = 1.180207364467719
$ws.Range("I22").Value = 0.8073460155169965
$ws.Range("J22").Value = 0.3473127891966072
$ws.Range("L22").Value = 1.509974492457275
# Row 23
$ws.Range("B23").Value = 2.477034439207671
$ws.Range("D23").Value = 0.1786377741548364
$ws.Range("E23").Value = 0.2333803519876199
$ws.Range("F23").Value = 1.747843181918014
$ws.Range("G23").Value = 1.382081963630071
$ws.Range("H23").Value = 1.177742811714467
$ws.Range("I23").Value = 0.8153210772986732
$ws.Range("J23").Value = 0.3397197410850197
$ws.Range("L23").Value = 1.466914442324537
# Row 24
$ws.Range("B24").Value = 2.230237810247274
$ws.Range("D24").Value = 0.1673550209746253
$ws.Range("E24").Value = 0.2153953305531005
$ws.Range("F24").Value = 1.740452682261875
$ws.Range("G24").Value = 1.340045288638549
$ws.Range("H24").Value = 1.170104902302626
$ws.Range("I24").Value = 0.8470315347268613
$ws.Range("J24").Value = 0.3109183289925284
$ws.Range("L24").Value = 1.303892543931283
# Row 25
$ws.Range("B25").Value = 1.964093883676014
$ws.Range("D25").Value = 0.1550698250651834
$ws.Range("E25").Value = 0.1959078214325132
$ws.Range("F25").Value = 1.738539480775529
$ws.Range("G25").Value = 1.299977162602488
$ws.Range("H25").Value = 1.165451879860967
$ws.Range("I25").Value = 0.8844332267518471
$ws.Range("J25").Value = 0.2797668624521066
$ws.Range("L25").Value = 1.128200127038582
